$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"35"
$ws.Cells.Item(2, 2).Value = '$\eta_{q}$'
$ws.Cells.Item(2, 3).Value = [double]"0.0003227950214445237"
$ws.Cells.Item(3, 1).Value = [double]"64"
$ws.Cells.Item(3, 2).Value = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Cells.Item(3, 3).Value = [double]"0.0001753855856154807"
$ws.Cells.Item(4, 1).Value = [double]"32"
$ws.Cells.Item(4, 2).Value = '$F_{q}^{\text{SCF}}$'
$ws.Cells.Item(4, 3).Value = [double]"0.0001673776427848547"
$ws.Cells.Item(5, 1).Value = [double]"26"
$ws.Cells.Item(5, 2).Value = '$type_3$'
$ws.Cells.Item(5, 3).Value = [double]"0.0001610445718908489"
$ws.Cells.Item(6, 1).Value = [double]"34"
$ws.Cells.Item(6, 2).Value = '$F_{q}$'
$ws.Cells.Item(6, 3).Value = [double]"7.570715940509855e-05"
$ws.Cells.Item(7, 1).Value = [double]"23"
$ws.Cells.Item(7, 2).Value = '$type_0$'
$ws.Cells.Item(7, 3).Value = [double]"6.879427968757793e-05"
$ws.Cells.Item(8, 1).Value = [double]"3"
$ws.Cells.Item(8, 2).Value = '(h$_{p}$)$_{3}$'
$ws.Cells.Item(8, 3).Value = [double]"5.77540952000027e-05"
$ws.Cells.Item(9, 1).Value = [double]"40"
$ws.Cells.Item(9, 2).Value = '$F_{s}^{\text{SCF}}$'
$ws.Cells.Item(9, 3).Value = [double]"5.263956527928243e-05"
$ws.Cells.Item(10, 1).Value = [double]"43"
$ws.Cells.Item(10, 2).Value = '$\eta_{s}$'
$ws.Cells.Item(10, 3).Value = [double]"4.926152601733272e-05"
$ws.Cells.Item(11, 1).Value = [double]"24"
$ws.Cells.Item(11, 2).Value = '$type_1$'
$ws.Cells.Item(11, 3).Value = [double]"4.701170733170275e-05"
$ws.Cells.Item(12, 1).Value = [double]"73"
$ws.Cells.Item(12, 2).Value = '$\langle ss \vert ss \rangle$'
$ws.Cells.Item(12, 3).Value = [double]"3.992743153664366e-05"
$ws.Cells.Item(13, 1).Value = [double]"12"
$ws.Cells.Item(13, 2).Value = 'h$_{q}$'
$ws.Cells.Item(13, 3).Value = [double]"2.815855409680222e-05"
$ws.Cells.Item(14, 1).Value = [double]"22"
$ws.Cells.Item(14, 2).Value = 'h$_{s}$'
$ws.Cells.Item(14, 3).Value = [double]"2.696122646959244e-05"
$ws.Cells.Item(15, 1).Value = [double]"89"
$ws.Cells.Item(15, 2).Value = '$(\langle rr \vert rr \rangle)_{2}$'
$ws.Cells.Item(15, 3).Value = [double]"2.448931413413242e-05"
$ws.Cells.Item(16, 1).Value = [double]"42"
$ws.Cells.Item(16, 2).Value = '$F_{s}$'
$ws.Cells.Item(16, 3).Value = [double]"2.416242784405762e-05"
$ws.Cells.Item(17, 1).Value = [double]"25"
$ws.Cells.Item(17, 2).Value = '$type_2$'
$ws.Cells.Item(17, 3).Value = [double]"2.413045251152342e-05"
$ws.Cells.Item(18, 1).Value = [double]"88"
$ws.Cells.Item(18, 2).Value = '$(\langle pp \vert pp \rangle)_{2}$'
$ws.Cells.Item(18, 3).Value = [double]"2.312587742922136e-05"
$ws.Cells.Item(19, 1).Value = [double]"71"
$ws.Cells.Item(19, 2).Value = '$\langle qq \vert qq \rangle$'
$ws.Cells.Item(19, 3).Value = [double]"2.047014961696811e-05"
$ws.Cells.Item(20, 1).Value = [double]"75"
$ws.Cells.Item(20, 2).Value = '$(\langle pq \vert qp \rangle)_{0}$'
$ws.Cells.Item(20, 3).Value = [double]"1.962903448433232e-05"
$ws.Cells.Item(21, 1).Value = [double]"90"
$ws.Cells.Item(21, 2).Value = '$(\langle pq \vert pq \rangle)_{2}$'
$ws.Cells.Item(21, 3).Value = [double]"1.725185565443219e-05"
$ws.Cells.Item(22, 1).Value = [double]"96"
$ws.Cells.Item(22, 2).Value = '$(\langle pp \vert pp \rangle)_{3}$'
$ws.Cells.Item(22, 3).Value = [double]"1.704692822557048e-05"
$ws.Cells.Item(23, 1).Value = [double]"0"
$ws.Cells.Item(23, 2).Value = '(h$_{p}$)$_{0}$'
$ws.Cells.Item(23, 3).Value = [double]"1.632258485374428e-05"
$ws.Cells.Item(24, 1).Value = [double]"48"
$ws.Cells.Item(24, 2).Value = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Cells.Item(24, 3).Value = [double]"1.457445281302858e-05"
$ws.Cells.Item(25, 1).Value = [double]"91"
$ws.Cells.Item(25, 2).Value = '$(\langle pq \vert qp \rangle)_{2}$'
$ws.Cells.Item(25, 3).Value = [double]"1.439184378747983e-05"
$ws.Cells.Item(26, 1).Value = [double]"92"
$ws.Cells.Item(26, 2).Value = '$(\langle rs\vert rs \rangle)_{2}$'
$ws.Cells.Item(26, 3).Value = [double]"1.25719996788816e-05"
$ws.Cells.Item(27, 1).Value = [double]"1"
$ws.Cells.Item(27, 2).Value = '(h$_{p}$)$_{1}$'
$ws.Cells.Item(27, 3).Value = [double]"1.114783795951332e-05"
$ws.Cells.Item(28, 1).Value = [double]"93"
$ws.Cells.Item(28, 2).Value = '$(\langle rs \vert sr \rangle)_{2}$'
$ws.Cells.Item(28, 3).Value = [double]"1.078555204045189e-05"
$ws.Cells.Item(29, 1).Value = [double]"58"
$ws.Cells.Item(29, 2).Value = '$(F_{r})_{2}$'
$ws.Cells.Item(29, 3).Value = [double]"1.071578629546989e-05"
$ws.Cells.Item(30, 1).Value = [double]"98"
$ws.Cells.Item(30, 2).Value = '$(\langle pq \vert pq \rangle)_{3}$'
$ws.Cells.Item(30, 3).Value = [double]"1.006839772539599e-05"
$ws.Cells.Item(31, 1).Value = [double]"56"
$ws.Cells.Item(31, 2).Value = '$(F_{r}^{\text{SCF}})_{2}$'
$ws.Cells.Item(31, 3).Value = [double]"9.769173504323023e-06"
$ws.Cells.Item(32, 1).Value = [double]"11"
$ws.Cells.Item(32, 2).Value = '(h$_{pr}$)$_{3}$'
$ws.Cells.Item(32, 3).Value = [double]"9.51062918574289e-06"
$ws.Cells.Item(33, 1).Value = [double]"15"
$ws.Cells.Item(33, 2).Value = '(h$_{r}$)$_{1}$'
$ws.Cells.Item(33, 3).Value = [double]"9.173416877423239e-06"
$ws.Cells.Item(34, 1).Value = [double]"83"
$ws.Cells.Item(34, 2).Value = '$(\langle pq \vert qp \rangle)_{1}$'
$ws.Cells.Item(34, 3).Value = [double]"8.678508376483038e-06"
$ws.Cells.Item(35, 1).Value = [double]"16"
$ws.Cells.Item(35, 2).Value = '(h$_{r}$)$_{2}$'
$ws.Cells.Item(35, 3).Value = [double]"7.896088743440199e-06"
$ws.Cells.Item(36, 1).Value = [double]"54"
$ws.Cells.Item(36, 2).Value = '$(F_{p})_{2}$'
$ws.Cells.Item(36, 3).Value = [double]"7.6667606184102e-06"
$ws.Cells.Item(37, 1).Value = [double]"67"
$ws.Cells.Item(37, 2).Value = '$(\eta_{r})_{3}$'
$ws.Cells.Item(37, 3).Value = [double]"7.176969880878323e-06"
$ws.Cells.Item(38, 1).Value = [double]"17"
$ws.Cells.Item(38, 2).Value = '(h$_{r}$)$_{3}$'
$ws.Cells.Item(38, 3).Value = [double]"7.110177598836302e-06"
$ws.Cells.Item(39, 1).Value = [double]"94"
$ws.Cells.Item(39, 2).Value = '$(\langle pq \vert rs \rangle)_{3}$'
$ws.Cells.Item(39, 3).Value = [double]"6.957688186202625e-06"
$ws.Cells.Item(40, 1).Value = [double]"99"
$ws.Cells.Item(40, 2).Value = '$(\langle pq \vert qp \rangle)_{3}$'
$ws.Cells.Item(40, 3).Value = [double]"6.935796939126503e-06"
$ws.Cells.Item(41, 1).Value = [double]"2"
$ws.Cells.Item(41, 2).Value = '(h$_{p}$)$_{2}$'
$ws.Cells.Item(41, 3).Value = [double]"6.729399788789961e-06"
$ws.Cells.Item(42, 1).Value = [double]"52"
$ws.Cells.Item(42, 2).Value = '$(F_{p}^{\text{SCF}})_{2}$'
$ws.Cells.Item(42, 3).Value = [double]"6.061597221617856e-06"
$ws.Cells.Item(43, 1).Value = [double]"13"
$ws.Cells.Item(43, 2).Value = 'h$_{qs}$'
$ws.Cells.Item(43, 3).Value = [double]"5.941356018362779e-06"
$ws.Cells.Item(44, 1).Value = [double]"10"
$ws.Cells.Item(44, 2).Value = '(h$_{pr}$)$_{2}$'
$ws.Cells.Item(44, 3).Value = [double]"5.891014227902865e-06"
$ws.Cells.Item(45, 1).Value = [double]"97"
$ws.Cells.Item(45, 2).Value = '$(\langle rr \vert rr \rangle)_{3}$'
$ws.Cells.Item(45, 3).Value = [double]"5.618362713144107e-06"
$ws.Cells.Item(46, 1).Value = [double]"59"
$ws.Cells.Item(46, 2).Value = '$(\eta_{r})_{2}$'
$ws.Cells.Item(46, 3).Value = [double]"5.338060986601502e-06"
$ws.Cells.Item(47, 1).Value = [double]"38"
$ws.Cells.Item(47, 2).Value = '$(F_{r})_{0}$'
$ws.Cells.Item(47, 3).Value = [double]"4.530380903826383e-06"
$ws.Cells.Item(48, 1).Value = [double]"44"
$ws.Cells.Item(48, 2).Value = '$(F_{p}^{\text{SCF}})_{1}$'
$ws.Cells.Item(48, 3).Value = [double]"4.302899165472879e-06"
$ws.Cells.Item(49, 1).Value = [double]"66"
$ws.Cells.Item(49, 2).Value = '$(F_{r})_{3}$'
$ws.Cells.Item(49, 3).Value = [double]"4.235280167744376e-06"
$ws.Cells.Item(50, 1).Value = [double]"74"
$ws.Cells.Item(50, 2).Value = '$(\langle pq \vert pq \rangle)_{0}$'
$ws.Cells.Item(50, 3).Value = [double]"4.069239602619167e-06"
$ws.Cells.Item(51, 1).Value = [double]"39"
$ws.Cells.Item(51, 2).Value = '$(\eta_{r})_{0}$'
$ws.Cells.Item(51, 3).Value = [double]"3.955465183942162e-06"
$ws.Cells.Item(52, 1).Value = [double]"36"
$ws.Cells.Item(52, 2).Value = '$(F_{r}^{\text{SCF}})_{0}$'
$ws.Cells.Item(52, 3).Value = [double]"3.931529687778767e-06"
$ws.Cells.Item(53, 1).Value = [double]"28"
$ws.Cells.Item(53, 2).Value = '$(F_{p}^{\text{SCF}})_{0}$'
$ws.Cells.Item(53, 3).Value = [double]"3.600087842877399e-06"
$ws.Cells.Item(54, 1).Value = [double]"14"
$ws.Cells.Item(54, 2).Value = '(h$_{r}$)$_{0}$'
$ws.Cells.Item(54, 3).Value = [double]"3.531127554336526e-06"
$ws.Cells.Item(55, 1).Value = [double]"46"
$ws.Cells.Item(55, 2).Value = '$(F_{p})_{1}$'
$ws.Cells.Item(55, 3).Value = [double]"3.48541636278326e-06"
$ws.Cells.Item(56, 1).Value = [double]"62"
$ws.Cells.Item(56, 2).Value = '$(F_{p})_{3}$'
$ws.Cells.Item(56, 3).Value = [double]"3.252432154811559e-06"
$ws.Cells.Item(57, 1).Value = [double]"8"
$ws.Cells.Item(57, 2).Value = '(h$_{pr}$)$_{0}$'
$ws.Cells.Item(57, 3).Value = [double]"3.202765399077644e-06"
$ws.Cells.Item(58, 1).Value = [double]"9"
$ws.Cells.Item(58, 2).Value = '(h$_{pr}$)$_{1}$'
$ws.Cells.Item(58, 3).Value = [double]"2.920699910651781e-06"
$ws.Cells.Item(59, 1).Value = [double]"70"
$ws.Cells.Item(59, 2).Value = '$(\langle pp \vert pp \rangle)_{0}$'
$ws.Cells.Item(59, 3).Value = [double]"2.88424697880187e-06"
$ws.Cells.Item(60, 1).Value = [double]"30"
$ws.Cells.Item(60, 2).Value = '$(F_{p})_{0}$'
$ws.Cells.Item(60, 3).Value = [double]"2.882693522547001e-06"
$ws.Cells.Item(61, 1).Value = [double]"5"
$ws.Cells.Item(61, 2).Value = '(h$_{pq}$)$_{1}$'
$ws.Cells.Item(61, 3).Value = [double]"2.837365505930205e-06"
$ws.Cells.Item(62, 1).Value = [double]"101"
$ws.Cells.Item(62, 2).Value = '$(\langle rs \vert sr \rangle)_{3}$'
$ws.Cells.Item(62, 3).Value = [double]"2.748105365620831e-06"
$ws.Cells.Item(63, 1).Value = [double]"60"
$ws.Cells.Item(63, 2).Value = '$(F_{p}^{\text{SCF}})_{3}$'
$ws.Cells.Item(63, 3).Value = [double]"2.69693990239924e-06"
$ws.Cells.Item(64, 1).Value = [double]"77"
$ws.Cells.Item(64, 2).Value = '$(\langle rs \vert sr \rangle)_{0}$'
$ws.Cells.Item(64, 3).Value = [double]"2.51334631288551e-06"
$ws.Cells.Item(65, 1).Value = [double]"72"
$ws.Cells.Item(65, 2).Value = '$(\langle rr \vert rr \rangle)_{0}$'
$ws.Cells.Item(65, 3).Value = [double]"2.439737469507662e-06"
$ws.Cells.Item(66, 1).Value = [double]"85"
$ws.Cells.Item(66, 2).Value = '$(\langle rs \vert sr \rangle)_{1}$'
$ws.Cells.Item(66, 3).Value = [double]"2.168225676520742e-06"
$ws.Cells.Item(67, 1).Value = [double]"19"
$ws.Cells.Item(67, 2).Value = '(h$_{rs}$)$_{1}$'
$ws.Cells.Item(67, 3).Value = [double]"2.090544731761845e-06"
$ws.Cells.Item(68, 1).Value = [double]"76"
$ws.Cells.Item(68, 2).Value = '$(\langle rs\vert rs \rangle)_{0}$'
$ws.Cells.Item(68, 3).Value = [double]"2.088416151256583e-06"
$ws.Cells.Item(69, 1).Value = [double]"50"
$ws.Cells.Item(69, 2).Value = '$(F_{r})_{1}$'
$ws.Cells.Item(69, 3).Value = [double]"1.928148437475984e-06"
$ws.Cells.Item(70, 1).Value = [double]"82"
$ws.Cells.Item(70, 2).Value = '$(\langle pq \vert pq \rangle)_{1}$'
$ws.Cells.Item(70, 3).Value = [double]"1.614953245947436e-06"
$ws.Cells.Item(71, 1).Value = [double]"80"
$ws.Cells.Item(71, 2).Value = '$(\langle pp \vert pp \rangle)_{1}$'
$ws.Cells.Item(71, 3).Value = [double]"1.534974687521078e-06"
$ws.Cells.Item(72, 1).Value = [double]"51"
$ws.Cells.Item(72, 2).Value = '$(\eta_{r})_{1}$'
$ws.Cells.Item(72, 3).Value = [double]"1.447370215696673e-06"
$ws.Cells.Item(73, 1).Value = [double]"81"
$ws.Cells.Item(73, 2).Value = '$(\langle rr \vert rr \rangle)_{1}$'
$ws.Cells.Item(73, 3).Value = [double]"1.384610267673221e-06"
$ws.Cells.Item(74, 1).Value = [double]"84"
$ws.Cells.Item(74, 2).Value = '$(\langle rs\vert rs \rangle)_{1}$'
$ws.Cells.Item(74, 3).Value = [double]"1.368442915242756e-06"
$ws.Cells.Item(75, 1).Value = [double]"21"
$ws.Cells.Item(75, 2).Value = '(h$_{rs}$)$_{3}$'
$ws.Cells.Item(75, 3).Value = [double]"1.288572415691375e-06"
$ws.Cells.Item(76, 1).Value = [double]"78"
$ws.Cells.Item(76, 2).Value = '$(\langle pq \vert rs \rangle)_{1}$'
$ws.Cells.Item(76, 3).Value = [double]"1.139929265136116e-06"
$ws.Cells.Item(77, 1).Value = [double]"18"
$ws.Cells.Item(77, 2).Value = '(h$_{rs}$)$_{0}$'
$ws.Cells.Item(77, 3).Value = [double]"8.960245586189157e-07"
$ws.Cells.Item(78, 1).Value = [double]"86"
$ws.Cells.Item(78, 2).Value = '$(\langle pq \vert rs \rangle)_{2}$'
$ws.Cells.Item(78, 3).Value = [double]"7.311489799564644e-07"
$ws.Cells.Item(79, 1).Value = [double]"100"
$ws.Cells.Item(79, 2).Value = '$(\langle rs\vert rs \rangle)_{3}$'
$ws.Cells.Item(79, 3).Value = [double]"6.791928301510475e-07"
$ws.Cells.Item(80, 1).Value = [double]"4"
$ws.Cells.Item(80, 2).Value = '(h$_{pq}$)$_{0}$'
$ws.Cells.Item(80, 3).Value = [double]"6.739016698251171e-07"
$ws.Cells.Item(81, 1).Value = [double]"7"
$ws.Cells.Item(81, 2).Value = '(h$_{pq}$)$_{3}$'
$ws.Cells.Item(81, 3).Value = [double]"5.453312792359722e-07"
$ws.Cells.Item(82, 1).Value = [double]"68"
$ws.Cells.Item(82, 2).Value = '$(\langle pq \vert rs \rangle)_{0}$'
$ws.Cells.Item(82, 3).Value = [double]"3.49406901722563e-07"
$ws.Cells.Item(83, 1).Value = [double]"20"
$ws.Cells.Item(83, 2).Value = '(h$_{rs}$)$_{2}$'
$ws.Cells.Item(83, 3).Value = [double]"2.834373844795528e-08"
$ws.Cells.Item(84, 1).Value = [double]"6"
$ws.Cells.Item(84, 2).Value = '(h$_{pq}$)$_{2}$'
$ws.Cells.Item(84, 3).Value = [double]"9.540340858167671e-09"
$ws.Cells.Item(85, 1).Value = [double]"95"
$ws.Cells.Item(85, 2).Value = '$(\langle pq \vert sr \rangle)_{3}$'
$ws.Cells.Item(85, 3).Value = [double]"3.521099726045723e-09"
$ws.Cells.Item(86, 1).Value = [double]"49"
$ws.Cells.Item(86, 2).Value = '$(\omega_{r})_{1}$'
$ws.Cells.Item(86, 3).Value = [double]"2.373144711014198e-09"
$ws.Cells.Item(87, 1).Value = [double]"57"
$ws.Cells.Item(87, 2).Value = '$(\omega_{r})_{2}$'
$ws.Cells.Item(87, 3).Value = [double]"2.306757980124273e-09"
$ws.Cells.Item(88, 1).Value = [double]"69"
$ws.Cells.Item(88, 2).Value = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Cells.Item(88, 3).Value = [double]"2.027877903381537e-09"
$ws.Cells.Item(89, 1).Value = [double]"41"
$ws.Cells.Item(89, 2).Value = '$\omega_{s}$'
$ws.Cells.Item(89, 3).Value = [double]"1.981271324415822e-09"
$ws.Cells.Item(90, 1).Value = [double]"37"
$ws.Cells.Item(90, 2).Value = '$(\omega_{r})_{0}$'
$ws.Cells.Item(90, 3).Value = [double]"1.427049858443511e-09"
$ws.Cells.Item(91, 1).Value = [double]"65"
$ws.Cells.Item(91, 2).Value = '$(\omega_{r})_{3}$'
$ws.Cells.Item(91, 3).Value = [double]"1.360869585401551e-09"
$ws.Cells.Item(92, 1).Value = [double]"87"
$ws.Cells.Item(92, 2).Value = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Cells.Item(92, 3).Value = [double]"1.272613058268239e-09"
$ws.Cells.Item(93, 1).Value = [double]"79"
$ws.Cells.Item(93, 2).Value = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Cells.Item(93, 3).Value = [double]"1.246339974398595e-09"
$ws.Cells.Item(94, 1).Value = [double]"33"
$ws.Cells.Item(94, 2).Value = '$\omega_{q}$'
$ws.Cells.Item(94, 3).Value = [double]"1.229068629982067e-09"
$ws.Cells.Item(95, 1).Value = [double]"27"
$ws.Cells.Item(95, 2).Value = '$\mathbf{b}$'
$ws.Cells.Item(95, 3).Value = [double]"1.7927719694614e-10"
